$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.2917716402565462, 0.04071648406533734, 0.7527432677738641, 10.19245300693656, 0, 11.2776843990323)
    3 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    4 = @(0.1190320826869504, 0.002571899574220771, 0.7527432677738641, 0.4942365360607697, 1, 1.368583786095805)
    5 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 3.56341032713086)
    6 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    7 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 6.348428708163715)
    8 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    9 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
